# Insert a new data row at row 23 (pushing existing rows 23..115 down to 24..116)
# and populate it with the new observation, matching the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 23; this shifts rows 23-115 down to 24-116
# and keeps all their values/styles intact, and also updates the sheet dimension.
$ws.Rows("23").Insert()

# Populate the newly inserted, now-empty row 23 with the new record's data.
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = 45107
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 100112043
$ws.Range("G23").Value = "Pepino dulce"
$ws.Range("H23").Value = "Cultivar IV Región"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 18000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 19000
$ws.Range("N23").Value = "`$/bandeja 18 kilos"
$ws.Range("O23").Value = "Provincia de Limarí"
$ws.Range("P23").Value = 1056
$ws.Range("Q23").Value = 18
$ws.Range("R23").Value = "Hortaliza"
